{"js": "// Map of old text -> new text, exactly as they appear in the document.\nconst replacements = [\n  [\"2025-06-24 Tuesday\", \"2025-06-25 Wednesday\"],\n  [\"294\u00d75=1470\", \"999\u00d79=8991\"],\n  [\"473\u00d76=2838\", \"906\u00d73=2718\"],\n  [\"401\u00d73=1203\", \"405\u00d79=3645\"],\n  [\"344\u00d75=1720\", \"743\u00d72=1486\"],\n  [\"808\u00d77=5656\", \"724\u00d74=2896\"],\n  [\"572\u00d76=3432\", \"254\u00d74=1016\"],\n  [\"686\u00d76=4116\", \"492\u00d78=3936\"],\n  [\"336\u00d79=3024\", \"810\u00d78=6480\"],\n  [\"424\u00d76=2544\", \"159\u00d78=1272\"],\n  [\"342\u00d74=1368\", \"317\u00d76=1902\"],\n  [\"549\u00d72=1098\", \"558\u00d72=1116\"],\n  [\"748\u00d77=5236\", \"947\u00d74=3788\"],\n  [\"414\u00d74=1656\", \"998\u00d74=3992\"],\n  [\"158\u00d76=948\", \"373\u00d72=746\"],\n  [\"246\u00d74=984\", \"184\u00d78=1472\"],\n  [\"815\u00d78=6520\", \"474\u00d77=3318\"],\n  [\"553\u00d74=2212\", \"762\u00d74=3048\"],\n  [\"800\u00d74=3200\", \"679\u00d73=2037\"],\n  [\"789\u00d75=3945\", \"488\u00d79=4392\"],\n  [\"532\u00d76=3192\", \"568\u00d78=4544\"],\n  [\"852\u00d76=5112\", \"499\u00d78=3992\"],\n  [\"949\u00d72=1898\", \"352\u00d77=2464\"],\n  [\"364\u00d78=2912\", \"793\u00d77=5551\"],\n  [\"408\u00d79=3672\", \"787\u00d72=1574\"],\n  [\"998\u00d77=6986\", \"457\u00d75=2285\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" cell value in the table with\n# the new values, matching the authored commit exactly. Each old string is\n# unique in the document, so a plain Find/Replace (Replace = wdReplaceAll)\n# per pair is safe and unambiguous.\n$pairs = @(\n    @{ Old = \"2025-06-24 Tuesday\"; New = \"2025-06-25 Wednesday\" },\n    @{ Old = \"294\u00d75=1470\"; New = \"999\u00d79=8991\" },\n    @{ Old = \"473\u00d76=2838\"; New = \"906\u00d73=2718\" },\n    @{ Old = \"401\u00d73=1203\"; New = \"405\u00d79=3645\" },\n    @{ Old = \"344\u00d75=1720\"; New = \"743\u00d72=1486\" },\n    @{ Old = \"808\u00d77=5656\"; New = \"724\u00d74=2896\" },\n    @{ Old = \"572\u00d76=3432\"; New = \"254\u00d74=1016\" },\n    @{ Old = \"686\u00d76=4116\"; New = \"492\u00d78=3936\" },\n    @{ Old = \"336\u00d79=3024\"; New = \"810\u00d78=6480\" },\n    @{ Old = \"424\u00d76=2544\"; New = \"159\u00d78=1272\" },\n    @{ Old = \"342\u00d74=1368\"; New = \"317\u00d76=1902\" },\n    @{ Old = \"549\u00d72=1098\"; New = \"558\u00d72=1116\" },\n    @{ Old = \"748\u00d77=5236\"; New = \"947\u00d74=3788\" },\n    @{ Old = \"414\u00d74=1656\"; New = \"998\u00d74=3992\" },\n    @{ Old = \"158\u00d76=948\";  New = \"373\u00d72=746\" },\n    @{ Old = \"246\u00d74=984\";  New = \"184\u00d78=1472\" },\n    @{ Old = \"815\u00d78=6520\"; New = \"474\u00d77=3318\" },\n    @{ Old = \"553\u00d74=2212\"; New = \"762\u00d74=3048\" },\n    @{ Old = \"800\u00d74=3200\"; New = \"679\u00d73=2037\" },\n    @{ Old = \"789\u00d75=3945\"; New = \"488\u00d79=4392\" },\n    @{ Old = \"532\u00d76=3192\"; New = \"568\u00d78=4544\" },\n    @{ Old = \"852\u00d76=5112\"; New = \"499\u00d78=3992\" },\n    @{ Old = \"949\u00d72=1898\"; New = \"352\u00d77=2464\" },\n    @{ Old = \"364\u00d78=2912\"; New = \"793\u00d77=5551\" },\n    @{ Old = \"408\u00d79=3672\"; New = \"787\u00d72=1574\" },\n    @{ Old = \"998\u00d77=6986\"; New = \"457\u00d75=2285\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($p.Old, $true, $false, $false, $false, $false, $true, 1, $false, $p.New, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($p.Old)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
